$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Recommandations")
$ws2 = $wb.Worksheets.Item("Top_YTD")

# Remove the 3 trailing rows from Recommandations (37 -> 34 rows)
$ws1.Rows.Item(37).EntireRow.Delete()
$ws1.Rows.Item(36).EntireRow.Delete()
$ws1.Rows.Item(35).EntireRow.Delete()

# Rewrite Recommandations rows 2-34 with the refreshed BRVM figures
$ws1.Range("A2").Value = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$ws1.Range("B2").Value = 0
$ws1.Range("C2").Value = 4
$ws1.Range("D2").Value = 646.08
$ws1.Range("E2").Value = 163.84
$ws1.Range("F2").Value = "🟡 Observer"
$ws1.Range("G2").Value = "➖ Neutre"

$ws1.Range("A3").Value = "BRVM - SERVICES FINANCIERS"
$ws1.Range("B3").Value = 0
$ws1.Range("C3").Value = 4
$ws1.Range("D3").Value = 582.82
$ws1.Range("E3").Value = 145.82
$ws1.Range("F3").Value = "🟡 Observer"
$ws1.Range("G3").Value = "➖ Neutre"

$ws1.Range("A4").Value = "BRVM-PRESTIGE"
$ws1.Range("B4").Value = 0
$ws1.Range("C4").Value = 4
$ws1.Range("D4").Value = 570.74
$ws1.Range("E4").Value = 142.44
$ws1.Range("F4").Value = "🟡 Observer"
$ws1.Range("G4").Value = "➖ Neutre"

$ws1.Range("A5").Value = "BRVM - INDUSTRIELS"
$ws1.Range("B5").Value = 0
$ws1.Range("C5").Value = 4
$ws1.Range("D5").Value = 530.07
$ws1.Range("E5").Value = 135.02
$ws1.Range("F5").Value = "🟡 Observer"
$ws1.Range("G5").Value = "➖ Neutre"

$ws1.Range("A6").Value = "BRVM - ENERGIE"
$ws1.Range("B6").Value = 0
$ws1.Range("C6").Value = 4
$ws1.Range("D6").Value = 444.26
$ws1.Range("E6").Value = 112.38
$ws1.Range("F6").Value = "🟡 Observer"
$ws1.Range("G6").Value = "➖ Neutre"

$ws1.Range("A7").Value = "BRVM - SERVICES PUBLICS"
$ws1.Range("B7").Value = 0
$ws1.Range("C7").Value = 4
$ws1.Range("D7").Value = 427.39
$ws1.Range("E7").Value = 107.51
$ws1.Range("F7").Value = "🟡 Observer"
$ws1.Range("G7").Value = "➖ Neutre"

$ws1.Range("A8").Value = "BRVM - TELECOMMUNICATIONS"
$ws1.Range("B8").Value = 0
$ws1.Range("C8").Value = 4
$ws1.Range("D8").Value = 372.34
$ws1.Range("E8").Value = 93.26
$ws1.Range("F8").Value = "🟡 Observer"
$ws1.Range("G8").Value = "➖ Neutre"

$ws1.Range("A9").Value = "BRVM-PRINCIPAL     (**)"
$ws1.Range("B9").Value = 0
$ws1.Range("C9").Value = 1
$ws1.Range("D9").Value = 216.13
$ws1.Range("E9").Value = 216.13
$ws1.Range("F9").Value = "🟡 Observer"
$ws1.Range("G9").Value = "➖ Neutre"

$ws1.Range("A10").Value = "BRVM - CONSOMMATION DE BASE    (**)"
$ws1.Range("B10").Value = 0
$ws1.Range("C10").Value = 1
$ws1.Range("D10").Value = 215.68
$ws1.Range("E10").Value = 215.68
$ws1.Range("F10").Value = "🟡 Observer"
$ws1.Range("G10").Value = "➖ Neutre"

$ws1.Range("A11").Value = "BRVM – COMPOSITE TOTAL RETURN    (**)"
$ws1.Range("B11").Value = 0
$ws1.Range("C11").Value = 1
$ws1.Range("D11").Value = 130.87
$ws1.Range("E11").Value = 130.87
$ws1.Range("F11").Value = "🟡 Observer"
$ws1.Range("G11").Value = "➖ Neutre"

$ws1.Range("A12").Value = "EVIOSYS PACKAGING SIEM CI (SEMC)"
$ws1.Range("B12").Value = 4
$ws1.Range("C12").Value = 0
$ws1.Range("D12").Value = 28.75
$ws1.Range("E12").Value = 7.03
$ws1.Range("F12").Value = "🟢 Achat"
$ws1.Range("G12").Value = "✅ Renforcer"

$ws1.Range("A13").Value = "VIVO ENERGY CI (SHEC)"
$ws1.Range("B13").Value = 2
$ws1.Range("C13").Value = 0
$ws1.Range("D13").Value = 9.69
$ws1.Range("E13").Value = 4.73
$ws1.Range("F13").Value = "🟡 Observer"
$ws1.Range("G13").Value = "➖ Neutre"

$ws1.Range("A14").Value = "SICOR CI (SICC)"
$ws1.Range("B14").Value = 1
$ws1.Range("C14").Value = 0
$ws1.Range("D14").Value = 7.4
$ws1.Range("E14").Value = 7.4
$ws1.Range("F14").Value = "🟡 Observer"
$ws1.Range("G14").Value = "➖ Neutre"

$ws1.Range("A15").Value = "UNILEVER CI (UNLC)"
$ws1.Range("B15").Value = 1
$ws1.Range("C15").Value = 0
$ws1.Range("D15").Value = 7.14
$ws1.Range("E15").Value = 7.14
$ws1.Range("F15").Value = "🟡 Observer"
$ws1.Range("G15").Value = "➖ Neutre"

$ws1.Range("A16").Value = "ORAGROUP TOGO (ORGT)"
$ws1.Range("B16").Value = 1
$ws1.Range("C16").Value = 0
$ws1.Range("D16").Value = 6.04
$ws1.Range("E16").Value = 6.04
$ws1.Range("F16").Value = "🟡 Observer"
$ws1.Range("G16").Value = "➖ Neutre"

$ws1.Range("A17").Value = "SOGB CI (SOGC)"
$ws1.Range("B17").Value = 1
$ws1.Range("C17").Value = 0
$ws1.Range("D17").Value = 3.24
$ws1.Range("E17").Value = 3.24
$ws1.Range("F17").Value = "🟡 Observer"
$ws1.Range("G17").Value = "➖ Neutre"

$ws1.Range("A18").Value = "FILTISAC CI (FTSC)"
$ws1.Range("B18").Value = 2
$ws1.Range("C18").Value = 2
$ws1.Range("D18").Value = 1.84
$ws1.Range("E18").Value = 7.26
$ws1.Range("F18").Value = "🟡 Observer"
$ws1.Range("G18").Value = "👀 À surveiller"

$ws1.Range("A19").Value = "SOLIBRA CI (SLBC)"
$ws1.Range("B19").Value = 2
$ws1.Range("C19").Value = 1
$ws1.Range("D19").Value = 1.46
$ws1.Range("E19").Value = 4.91
$ws1.Range("F19").Value = "🟡 Observer"
$ws1.Range("G19").Value = "👀 À surveiller"

$ws1.Range("A20").Value = "SOCIETE IVOIRIENNE DE BANQUE  (SIBC)"
$ws1.Range("B20").Value = 1
$ws1.Range("C20").Value = 1
$ws1.Range("D20").Value = 1.42
$ws1.Range("E20").Value = -2.68
$ws1.Range("F20").Value = "🟡 Observer"
$ws1.Range("G20").Value = "👀 À surveiller"

$ws1.Range("A21").Value = "ONATEL BF (ONTBF)"
$ws1.Range("B21").Value = 1
$ws1.Range("C21").Value = 0
$ws1.Range("D21").Value = 1.22
$ws1.Range("E21").Value = 1.22
$ws1.Range("F21").Value = "🟡 Observer"
$ws1.Range("G21").Value = "➖ Neutre"

$ws1.Range("A22").Value = "TOTALENERGIES MARKETING CI (TTLC)"
$ws1.Range("B22").Value = 0
$ws1.Range("C22").Value = 1
$ws1.Range("D22").Value = -1.49
$ws1.Range("E22").Value = -1.49
$ws1.Range("F22").Value = "🟡 Observer"
$ws1.Range("G22").Value = "➖ Neutre"

$ws1.Range("A23").Value = "SERVAIR ABIDJAN CI (ABJC)"
$ws1.Range("B23").Value = 0
$ws1.Range("C23").Value = 1
$ws1.Range("D23").Value = -1.8
$ws1.Range("E23").Value = -1.8
$ws1.Range("F23").Value = "🟡 Observer"
$ws1.Range("G23").Value = "➖ Neutre"

$ws1.Range("A24").Value = "UNIWAX CI (UNXC)"
$ws1.Range("B24").Value = 0
$ws1.Range("C24").Value = 1
$ws1.Range("D24").Value = -2.13
$ws1.Range("E24").Value = -2.13
$ws1.Range("F24").Value = "🟡 Observer"
$ws1.Range("G24").Value = "➖ Neutre"

$ws1.Range("A25").Value = "LOTERIE NATIONALE DU BENIN (LNBB)"
$ws1.Range("B25").Value = 0
$ws1.Range("C25").Value = 1
$ws1.Range("D25").Value = -2.21
$ws1.Range("E25").Value = -2.21
$ws1.Range("F25").Value = "🟡 Observer"
$ws1.Range("G25").Value = "➖ Neutre"

$ws1.Range("A26").Value = "SICABLE CI (CABC)"
$ws1.Range("B26").Value = 0
$ws1.Range("C26").Value = 1
$ws1.Range("D26").Value = -2.34
$ws1.Range("E26").Value = -2.34
$ws1.Range("F26").Value = "🟡 Observer"
$ws1.Range("G26").Value = "➖ Neutre"

$ws1.Range("A27").Value = "NEI-CEDA CI (NEIC)"
$ws1.Range("B27").Value = 1
$ws1.Range("C27").Value = 2
$ws1.Range("D27").Value = -2.66
$ws1.Range("E27").Value = 7.46
$ws1.Range("F27").Value = "🟡 Observer"
$ws1.Range("G27").Value = "👀 À surveiller"

$ws1.Range("A28").Value = "SOCIETE GENERALE COTE D'IVOIRE (SGBC)"
$ws1.Range("B28").Value = 0
$ws1.Range("C28").Value = 1
$ws1.Range("D28").Value = -3.06
$ws1.Range("E28").Value = -3.06
$ws1.Range("F28").Value = "🟡 Observer"
$ws1.Range("G28").Value = "➖ Neutre"

$ws1.Range("A29").Value = "SMB CI (SMBC)"
$ws1.Range("B29").Value = 0
$ws1.Range("C29").Value = 1
$ws1.Range("D29").Value = -3.16
$ws1.Range("E29").Value = -3.16
$ws1.Range("F29").Value = "🟡 Observer"
$ws1.Range("G29").Value = "➖ Neutre"

$ws1.Range("A30").Value = "ERIUM CI (Ex AIR LIQUIDE CI) (SIVC)"
$ws1.Range("B30").Value = 1
$ws1.Range("C30").Value = 2
$ws1.Range("D30").Value = -3.59
$ws1.Range("E30").Value = -3.9
$ws1.Range("F30").Value = "🟡 Observer"
$ws1.Range("G30").Value = "👀 À surveiller"

$ws1.Range("A31").Value = "TOTALENERGIES MARKETING SN (TTLS)"
$ws1.Range("B31").Value = 0
$ws1.Range("C31").Value = 1
$ws1.Range("D31").Value = -3.92
$ws1.Range("E31").Value = 3.19
$ws1.Range("F31").Value = "🟡 Observer"
$ws1.Range("G31").Value = "👀 À surveiller"

$ws1.Range("A32").Value = "ECOBANK TRANS. INCORP. TG (ETIT)"
$ws1.Range("B32").Value = 0
$ws1.Range("C32").Value = 1
$ws1.Range("D32").Value = -4.35
$ws1.Range("E32").Value = -4.35
$ws1.Range("F32").Value = "🟡 Observer"
$ws1.Range("G32").Value = "➖ Neutre"

$ws1.Range("A33").Value = "SETAO CI (STAC)"
$ws1.Range("B33").Value = 0
$ws1.Range("C33").Value = 1
$ws1.Range("D33").Value = -5.49
$ws1.Range("E33").Value = -5.49
$ws1.Range("F33").Value = "🟡 Observer"
$ws1.Range("G33").Value = "➖ Neutre"

$ws1.Range("A34").Value = "CFAO MOTORS CI (CFAC)"
$ws1.Range("B34").Value = 1
$ws1.Range("C34").Value = 2
$ws1.Range("D34").Value = -5.95
$ws1.Range("E34").Value = 4.71
$ws1.Range("F34").Value = "🟡 Observer"
$ws1.Range("G34").Value = "👀 À surveiller"

# Rewrite Top_YTD rows 2-11 with the refreshed YTD progression + label shift
$ws2.Range("A2").Value = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$ws2.Range("B2").Value = 4576.96

$ws2.Range("A3").Value = "BRVM - SERVICES FINANCIERS"
$ws2.Range("B3").Value = 3544.6

$ws2.Range("A4").Value = "BRVM-PRESTIGE"
$ws2.Range("B4").Value = 3368.65

$ws2.Range("A5").Value = "BRVM - INDUSTRIELS"
$ws2.Range("B5").Value = 2822.56

$ws2.Range("A6").Value = "BRVM - ENERGIE"
$ws2.Range("B6").Value = 1884.41

$ws2.Range("A7").Value = "BRVM - SERVICES PUBLICS"
$ws2.Range("B7").Value = 1730.6

$ws2.Range("A8").Value = "BRVM - TELECOMMUNICATIONS"
$ws2.Range("B8").Value = 1289.91

$ws2.Range("A9").Value = "BRVM-PRINCIPAL     (**)"
$ws2.Range("B9").Value = 216.13

$ws2.Range("A10").Value = "BRVM - CONSOMMATION DE BASE    (**)"
$ws2.Range("B10").Value = 215.68

$ws2.Range("A11").Value = "BRVM – COMPOSITE TOTAL RETURN    (**)"
$ws2.Range("B11").Value = 130.87

